# Update the SeenRx_CNH data matrix (Sheet1, A1:J46) with refreshed counts
# per the "finalize version to deliver reports to user" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> { column number (B=2 .. J=10) -> new value }
$updates = @{
    2 = @{ 2=79; 3=50; 4=56; 5=17; 6=23; 7=62; 8=105; 9=21; 10=61 }
    3 = @{ 2=23; 3=21; 5=8; 6=4; 7=19; 8=40; 9=8; 10=17 }
    4 = @{ 3=1; 5=2; 6=1; 7=1; 8=5; 9=2 }
    5 = @{ 2=1; 3=2; 4=3; 5=1; 7=1; 8=5; 10=5 }
    6 = @{ 2=3; 5=2; 7=7; 9=2; 10=1 }
    7 = @{ 2=6; 3=7; 4=4; 5=2; 7=4; 8=7 }
    8 = @{ 2=5; 3=8; 4=2; 6=3; 8=5; 9=3; 10=7 }
    9 = @{ 2=4; 7=4; 8=9 }
    10 = @{ 2=4; 3=2; 4=3; 5=1; 6=0; 8=6; 10=0 }
    11 = @{ 2=16; 3=15; 4=20; 6=4; 7=16; 9=3; 10=10 }
    12 = @{ 2=4; 3=1; 4=4; 8=0; 10=3 }
    13 = @{ 2=2; 4=1; 6=1; 8=0 }
    14 = @{ 2=1; 3=1; 4=2; 7=3; 8=3; 9=2; 10=2 }
    15 = @{ 2=5; 3=5; 4=5; 6=2; 7=5; 10=2 }
    16 = @{ 3=7; 4=3; 7=4; 10=2 }
    17 = @{ 2=4; 3=1; 4=5; 6=1; 7=2; 9=1; 10=1 }
    19 = @{ 2=11; 3=2; 4=2; 5=4; 6=2; 7=8; 8=5; 10=8 }
    20 = @{ 2=7; 5=1; 8=1; 10=0 }
    21 = @{ 4=1; 7=0; 8=1; 10=5 }
    22 = @{ 2=4; 3=1; 4=1; 7=5; 8=1; 10=1 }
    23 = @{ 3=1; 7=2; 8=2 }
    24 = @{ 5=3; 6=1; 7=1; 10=2 }
    25 = @{ 2=8; 4=5; 5=2; 6=2; 8=29; 9=1; 10=5 }
    26 = @{ 2=7; 4=2; 5=2; 6=2; 8=5; 10=1 }
    27 = @{ 8=5; 10=3 }
    29 = @{ 2=1; 8=1; 10=1 }
    30 = @{ 4=3; 8=18; 9=1; 10=0 }
    32 = @{ 2=13; 3=9; 4=13; 5=3; 6=9; 7=15; 8=18; 9=1; 10=10 }
    35 = @{ 2=5; 3=2; 4=1; 5=2; 6=3; 7=3 }
    36 = @{ 2=1; 3=2; 4=6; 6=2; 7=3; 8=6; 10=4 }
    37 = @{ 2=1; 3=2; 4=2; 8=3 }
    38 = @{ 3=2; 4=2; 5=1; 6=1; 7=2; 8=4; 9=1; 10=2 }
    39 = @{ 2=5; 3=1; 4=2; 6=2; 7=6; 8=5; 10=3 }
    40 = @{ 2=8; 3=3; 6=2; 7=4; 8=10; 9=8; 10=11 }
    43 = @{ 2=1; 6=1; 7=0; 8=4; 10=4 }
    44 = @{ 2=2; 6=1; 7=4; 8=5; 9=2; 10=4 }
    45 = @{ 2=3; 3=1; 8=1; 9=4; 10=2 }
    46 = @{ 3=1; 9=2; 10=1 }
}

foreach ($rowNum in $updates.Keys) {
    $rowMap = $updates[$rowNum]
    foreach ($colNum in $rowMap.Keys) {
        $ws.Cells.Item([int]$rowNum, [int]$colNum).Value = $rowMap[$colNum]
    }
}

Write-Output "Updated $($updates.Count) rows in $($ws.Name)"
